$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Title heading and bold "Play Cash Spin..." text (appears twice)
Replace-Text "Play Cash Spin Slot Game for Free - Exciting Bonuses and Classic Graphics" "Play Cash Spin Slot Game for Free"
Replace-Text "Play Cash Spin Slot Game for Free - Exciting Bonuses and Classic Graphics" "Play Cash Spin Slot Game for Free"

# What we like bullets
Replace-Text "Exciting bonus features with cash prizes and free spins" "Seamlessly blends classic slot themes with modern features"
Replace-Text "Classic yet appealing graphics with modern features" "Offers two bonus modes with cash prizes and free spins"
Replace-Text "Developed by a well-known and experienced game developer" "Includes a Wild symbol to increase winning potential"
Replace-Text "In line with modern standards for RTP" "Developed by a well-known and reputable developer, Bally"

# What we don't like bullets
Replace-Text "Limited theme and symbols may not appeal to all players" "Limited variety of symbols"
Replace-Text "The RTP is slightly lower than some other modern slot games" "Graphics may not appeal to players looking for more visually stimulating games"

# Meta description italic text
Replace-Text "Get a taste of classic slot gaming with modern features. Play Cash Spin for free and activate its exciting bonuses with cash prizes and free spins." "Read our review of Cash Spin, a slot game that seamlessly blends classic themes with modern features. Play for free!"
